# Auto-generated script to update market-price / leve-profit values
# in the "Louisoix_Profits" workbook, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC (88 cell updates) -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 772.1429000000001
$ws.Range("I18").Value = 734.1667
$ws.Range("K18").Value = 734.1667
$ws.Range("M18").Value = -450.1667
$ws.Range("H28").Value = 983.4286
$ws.Range("I28").Value = 979.6667
$ws.Range("J28").Value = 1006
$ws.Range("K28").Value = 979.6667
$ws.Range("L28").Value = 1006
$ws.Range("M28").Value = -494.6667
$ws.Range("N28").Value = -1976
$ws.Range("H34").Value = 8888.833000000001
$ws.Range("I34").Value = 6060.636
$ws.Range("K34").Value = 6060.636
$ws.Range("M34").Value = -5857.636
$ws.Range("H36").Value = 8888.833000000001
$ws.Range("I36").Value = 6060.636
$ws.Range("K36").Value = 6060.636
$ws.Range("M36").Value = -5345.636
$ws.Range("H38").Value = 331.25
$ws.Range("I38").Value = 331.25
$ws.Range("K38").Value = 993.75
$ws.Range("M38").Value = -621.75
$ws.Range("H39").Value = 1818.625
$ws.Range("I39").Value = 264
$ws.Range("K39").Value = 792
$ws.Range("M39").Value = -496
$ws.Range("H42").Value = 3211
$ws.Range("I42").Value = 792
$ws.Range("J42").Value = 5226.8335
$ws.Range("K42").Value = 2376
$ws.Range("L42").Value = 15680.5005
$ws.Range("M42").Value = -2146
$ws.Range("N42").Value = -16140.5005
$ws.Range("H51").Value = 10799.8
$ws.Range("I51").Value = 10666.667
$ws.Range("J51").Value = 10999.5
$ws.Range("K51").Value = 10666.667
$ws.Range("L51").Value = 10999.5
$ws.Range("M51").Value = -10182.667
$ws.Range("N51").Value = -11967.5
$ws.Range("H53").Value = 673.8
$ws.Range("I53").Value = 616.4286
$ws.Range("J53").Value = 724
$ws.Range("K53").Value = 616.4286
$ws.Range("L53").Value = 724
$ws.Range("M53").Value = 20.57140000000004
$ws.Range("N53").Value = -1998
$ws.Range("H76").Value = 4639
$ws.Range("I76").Value = 4298.75
$ws.Range("K76").Value = 4298.75
$ws.Range("M76").Value = -3983.75
$ws.Range("H79").Value = 4639
$ws.Range("I79").Value = 4298.75
$ws.Range("K79").Value = 4298.75
$ws.Range("M79").Value = -3206.75
$ws.Range("H86").Value = 113580
$ws.Range("I86").Value = 127502.5
$ws.Range("J86").Value = 2200
$ws.Range("K86").Value = 127502.5
$ws.Range("L86").Value = 2200
$ws.Range("M86").Value = -126379.5
$ws.Range("N86").Value = -4446
$ws.Range("H89").Value = 113580
$ws.Range("I89").Value = 127502.5
$ws.Range("J89").Value = 2200
$ws.Range("K89").Value = 637512.5
$ws.Range("L89").Value = 11000
$ws.Range("M89").Value = -631896.5
$ws.Range("N89").Value = -22232
$ws.Range("H107").Value = 1214.7241
$ws.Range("J107").Value = 1789.7778
$ws.Range("L107").Value = 1789.7778
$ws.Range("N107").Value = -5629.7778
$ws.Range("H135").Value = 1998.8636
$ws.Range("I135").Value = 1398.5625
$ws.Range("J135").Value = 3599.6667
$ws.Range("K135").Value = 12587.0625
$ws.Range("L135").Value = 32397.0003
$ws.Range("M135").Value = -10052.0625
$ws.Range("N135").Value = -37467.0003
$ws.Range("H138").Value = 5516.8364
$ws.Range("I138").Value = 3698.75
$ws.Range("J138").Value = 6555.7427
$ws.Range("K138").Value = 11096.25
$ws.Range("L138").Value = 19667.2281
$ws.Range("M138").Value = -5956.25
$ws.Range("N138").Value = -29947.2281

# ----- Sheet: ARM (32 cell updates) -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23125.25
$ws.Range("I32").Value = 24110.9
$ws.Range("J32").Value = 18197
$ws.Range("K32").Value = 24110.9
$ws.Range("L32").Value = 18197
$ws.Range("M32").Value = -23823.9
$ws.Range("N32").Value = -18771
$ws.Range("H62").Value = 13422.5
$ws.Range("I62").Value = 1845
$ws.Range("J62").Value = 25000
$ws.Range("K62").Value = 1845
$ws.Range("L62").Value = 25000
$ws.Range("M62").Value = -1221
$ws.Range("N62").Value = -26248
$ws.Range("H65").Value = 13422.5
$ws.Range("I65").Value = 1845
$ws.Range("J65").Value = 25000
$ws.Range("K65").Value = 5535
$ws.Range("L65").Value = 75000
$ws.Range("M65").Value = -2415
$ws.Range("N65").Value = -81240
$ws.Range("H102").Value = 12682.8
$ws.Range("I102").Value = 16429.143
$ws.Range("K102").Value = 16429.143
$ws.Range("M102").Value = -14807.143
$ws.Range("H132").Value = 40489.406
$ws.Range("I132").Value = 49949.145
$ws.Range("J132").Value = 7380.3335
$ws.Range("K132").Value = 149847.435
$ws.Range("L132").Value = 22141.0005
$ws.Range("M132").Value = -147317.435
$ws.Range("N132").Value = -27201.0005

# ----- Sheet: BSM (20 cell updates) -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 6315.8184
$ws.Range("I36").Value = 6315.8184
$ws.Range("K36").Value = 6315.8184
$ws.Range("M36").Value = -5781.8184
$ws.Range("H62").Value = 90064
$ws.Range("I62").Value = 90064
$ws.Range("K62").Value = 90064
$ws.Range("M62").Value = -89378
$ws.Range("H63").Value = 271000
$ws.Range("J63").Value = 271000
$ws.Range("L63").Value = 271000
$ws.Range("N63").Value = -272372
$ws.Range("H65").Value = 90064
$ws.Range("I65").Value = 90064
$ws.Range("K65").Value = 270192
$ws.Range("M65").Value = -266760
$ws.Range("H66").Value = 271000
$ws.Range("J66").Value = 271000
$ws.Range("L66").Value = 813000
$ws.Range("N66").Value = -819864

# ----- Sheet: CRP (34 cell updates) -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1342.2354
$ws.Range("I22").Value = 309.85715
$ws.Range("J22").Value = 2064.9
$ws.Range("K22").Value = 309.85715
$ws.Range("L22").Value = 2064.9
$ws.Range("M22").Value = 40.14285000000001
$ws.Range("N22").Value = -2764.9
$ws.Range("H31").Value = 3933.48
$ws.Range("I31").Value = 2036.4117
$ws.Range("K31").Value = 2036.4117
$ws.Range("M31").Value = -1741.4117
$ws.Range("H34").Value = 3933.48
$ws.Range("I34").Value = 2036.4117
$ws.Range("K34").Value = 2036.4117
$ws.Range("M34").Value = -1834.4117
$ws.Range("H62").Value = 3933.3333
$ws.Range("I62").Value = 3800
$ws.Range("K62").Value = 3800
$ws.Range("M62").Value = -3176
$ws.Range("H65").Value = 3933.3333
$ws.Range("I65").Value = 3800
$ws.Range("K65").Value = 19000
$ws.Range("M65").Value = -15880
$ws.Range("H132").Value = 4596.6
$ws.Range("I132").Value = 4329.5557
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 12988.6671
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -10458.6671
$ws.Range("N132").Value = -26060
$ws.Range("H134").Value = 93334.37
$ws.Range("I134").Value = 93334.37
$ws.Range("K134").Value = 280003.11
$ws.Range("M134").Value = -277468.11

# ----- Sheet: CUL (83 cell updates) -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 130.45454
$ws.Range("I2").Value = 157.22223
$ws.Range("K2").Value = 943.33338
$ws.Range("M2").Value = -830.33338
$ws.Range("H5").Value = 922.8889
$ws.Range("I5").Value = 1009.6
$ws.Range("K5").Value = 3028.8
$ws.Range("M5").Value = -2916.8
$ws.Range("H68").Value = 490.4
$ws.Range("I68").Value = 490.8889
$ws.Range("J68").Value = 486
$ws.Range("K68").Value = 1472.6667
$ws.Range("L68").Value = 1458
$ws.Range("M68").Value = -661.6667
$ws.Range("N68").Value = -3080
$ws.Range("H69").Value = 987.25
$ws.Range("J69").Value = 999.7143
$ws.Range("L69").Value = 2999.1429
$ws.Range("N69").Value = -4621.1429
$ws.Range("H71").Value = 490.4
$ws.Range("I71").Value = 490.8889
$ws.Range("J71").Value = 486
$ws.Range("K71").Value = 4418.0001
$ws.Range("L71").Value = 4374
$ws.Range("M71").Value = -362.0001000000002
$ws.Range("N71").Value = -12486
$ws.Range("H72").Value = 987.25
$ws.Range("J72").Value = 999.7143
$ws.Range("L72").Value = 8997.4287
$ws.Range("N72").Value = -17109.4287
$ws.Range("H75").Value = 698
$ws.Range("J75").Value = 698
$ws.Range("L75").Value = 2094
$ws.Range("N75").Value = -4090
$ws.Range("H76").Value = 11506.5
$ws.Range("I76").Value = 3013
$ws.Range("K76").Value = 9039
$ws.Range("M76").Value = -8656
$ws.Range("H78").Value = 698
$ws.Range("J78").Value = 698
$ws.Range("L78").Value = 6282
$ws.Range("N78").Value = -16266
$ws.Range("H79").Value = 11506.5
$ws.Range("I79").Value = 3013
$ws.Range("K79").Value = 9039
$ws.Range("M79").Value = -7713
$ws.Range("H96").Value = 18341.334
$ws.Range("I96").Value = 10025
$ws.Range("J96").Value = 22499.5
$ws.Range("K96").Value = 30075
$ws.Range("L96").Value = 67498.5
$ws.Range("M96").Value = -28016
$ws.Range("N96").Value = -71616.5
$ws.Range("H102").Value = 8099.5
$ws.Range("I102").Value = 6499
$ws.Range("J102").Value = 8633
$ws.Range("K102").Value = 19497
$ws.Range("L102").Value = 25899
$ws.Range("M102").Value = -17063
$ws.Range("N102").Value = -30767
$ws.Range("H104").Value = 3195.8333
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 3195.8333
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 9587.499899999999
$ws.Range("M104").Value = ""
$ws.Range("N104").Value = -14829.4999
$ws.Range("H123").Value = 2164.5
$ws.Range("I123").Value = 2164.5
$ws.Range("K123").Value = 6493.5
$ws.Range("M123").Value = -4043.5
$ws.Range("H131").Value = 5274193
$ws.Range("J131").Value = 9107931
$ws.Range("L131").Value = 27323793
$ws.Range("N131").Value = -27333873
$ws.Range("H135").Value = 922.8889
$ws.Range("I135").Value = 1009.6
$ws.Range("K135").Value = 9086.4
$ws.Range("M135").Value = -6551.4
$ws.Range("H139").Value = 1347.1111
$ws.Range("I139").Value = 1347.1111
$ws.Range("K139").Value = 4041.3333
$ws.Range("M139").Value = 1098.6667

# ----- Sheet: GSM (27 cell updates) -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 145
$ws.Range("I11").Value = 145
$ws.Range("K11").Value = 145
$ws.Range("M11").Value = -6
$ws.Range("H12").Value = 10000
$ws.Range("J12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("N12").Value = -10280
$ws.Range("H64").Value = 167050
$ws.Range("J64").Value = 167050
$ws.Range("L64").Value = 167050
$ws.Range("N64").Value = -167546
$ws.Range("H67").Value = 167050
$ws.Range("J67").Value = 167050
$ws.Range("L67").Value = 167050
$ws.Range("N67").Value = -168766
$ws.Range("H123").Value = 77777
$ws.Range("J123").Value = 77777
$ws.Range("L123").Value = 77777
$ws.Range("N123").Value = -82677
$ws.Range("H132").Value = 63136.117
$ws.Range("I132").Value = 85608.414
$ws.Range("J132").Value = 9202.6
$ws.Range("K132").Value = 256825.242
$ws.Range("L132").Value = 27607.8
$ws.Range("M132").Value = -254295.242
$ws.Range("N132").Value = -32667.8

# ----- Sheet: LTW (40 cell updates) -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1479.4062
$ws.Range("I22").Value = 689.5454999999999
$ws.Range("J22").Value = 3217.1
$ws.Range("K22").Value = 689.5454999999999
$ws.Range("L22").Value = 3217.1
$ws.Range("M22").Value = -394.5454999999999
$ws.Range("N22").Value = -3807.1
$ws.Range("H27").Value = 1479.4062
$ws.Range("I27").Value = 689.5454999999999
$ws.Range("J27").Value = 3217.1
$ws.Range("K27").Value = 689.5454999999999
$ws.Range("L27").Value = 3217.1
$ws.Range("M27").Value = -582.5454999999999
$ws.Range("N27").Value = -3431.1
$ws.Range("H64").Value = 24967.334
$ws.Range("J64").Value = 24967.334
$ws.Range("L64").Value = 24967.334
$ws.Range("N64").Value = -25417.334
$ws.Range("H67").Value = 24967.334
$ws.Range("J67").Value = 24967.334
$ws.Range("L67").Value = 24967.334
$ws.Range("N67").Value = -26527.334
$ws.Range("H93").Value = 1610.125
$ws.Range("I93").Value = 841.63635
$ws.Range("K93").Value = 841.63635
$ws.Range("M93").Value = 406.36365
$ws.Range("H122").Value = 4031.675
$ws.Range("I122").Value = 3417.3333
$ws.Range("J122").Value = 4534.3184
$ws.Range("K122").Value = 10251.9999
$ws.Range("L122").Value = 13602.9552
$ws.Range("M122").Value = -7801.999899999999
$ws.Range("N122").Value = -18502.9552
$ws.Range("H132").Value = 44233.586
$ws.Range("I132").Value = 56979.453
$ws.Range("J132").Value = 4175.143
$ws.Range("K132").Value = 170938.359
$ws.Range("L132").Value = 12525.429
$ws.Range("M132").Value = -168408.359
$ws.Range("N132").Value = -17585.429

# ----- Sheet: WVR (37 cell updates) -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = ""
$ws.Range("H63").Value = 17450
$ws.Range("J63").Value = 17450
$ws.Range("L63").Value = 17450
$ws.Range("N63").Value = -18698
$ws.Range("H66").Value = 17450
$ws.Range("J66").Value = 17450
$ws.Range("L66").Value = 52350
$ws.Range("N66").Value = -58590
$ws.Range("H81").Value = 2106.8333
$ws.Range("I81").Value = 975.0909
$ws.Range("J81").Value = 3885.2856
$ws.Range("K81").Value = 1950.1818
$ws.Range("L81").Value = 7770.5712
$ws.Range("M81").Value = -889.1818000000001
$ws.Range("N81").Value = -9892.5712
$ws.Range("H84").Value = 2106.8333
$ws.Range("I84").Value = 975.0909
$ws.Range("J84").Value = 3885.2856
$ws.Range("K84").Value = 9750.909
$ws.Range("L84").Value = 38852.856
$ws.Range("M84").Value = -4446.909
$ws.Range("N84").Value = -49460.856
$ws.Range("H122").Value = 1900
$ws.Range("I122").Value = 1900
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5700
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3250
$ws.Range("N122").Value = ""
$ws.Range("H132").Value = 145809.83
$ws.Range("I132").Value = 158707
$ws.Range("K132").Value = 476121
$ws.Range("M132").Value = -473591

Write-Host "Applied 361 cell updates across 8 sheets."